$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 35, shifting existing rows 35..140 down to 36..141
$ws.Rows.Item(35).Insert()

# Fill in the new row 35 with the same "template" values as the surrounding rows
# (these columns do not change across rows for this product/sheet)
$ws.Range("A35").Value = 10
$ws.Range("B35").Value = "Vega Modelo de Temuco"
$ws.Range("C35").Value = "La Araucanía"
$ws.Range("D35").Value = 44487
$ws.Range("E35").Value = 9
$ws.Range("F35").Value = 100112005
$ws.Range("G35").Value = "Puerro"
$ws.Range("H35").Value = "Azul de Maquehue"
$ws.Range("I35").Value = "Primera"
$ws.Range("J35").Value = 75
$ws.Range("K35").Value = 6000
$ws.Range("L35").Value = 7000
$ws.Range("M35").Value = 6600
$ws.Range("N35").Value = "`$/docena de paquetes"
$ws.Range("O35").Value = "Provincia de Cautín"
$ws.Range("P35").Value = 550
$ws.Range("Q35").Value = 12
$ws.Range("R35").Value = "Hortaliza"
